$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.028.39'
$ws.Range('E2').Value = '  +1.02%  '

# Row 3
$ws.Range('D3').Value = '1.636.50'
$ws.Range('E3').Value = '  +0.17%  '

# Row 4
$ws.Range('D4').Value = "'0.993"
$ws.Range('E4').Value = '  -0.77%  '

# Row 5
$ws.Range('D5').Value = "'214.97"
$ws.Range('E5').Value = '  -0.03%  '

# Row 6
$ws.Range('D6').Value = "'0.503"
$ws.Range('E6').Value = '  -0.31%  '

# Row 7
$ws.Range('D7').Value = "'0.995"
$ws.Range('E7').Value = '  -0.65%  '

# Row 8
$ws.Range('D8').Value = "'0.257"
$ws.Range('E8').Value = '  -0.66%  '

# Row 9
$ws.Range('D9').Value = "'0.0631"
$ws.Range('E9').Value = '  -0.89%  '

# Row 10
$ws.Range('D10').Value = "'19.70"
$ws.Range('E10').Value = '  +0.10%  '

# Row 11
$ws.Range('D11').Value = "'0.0787"
$ws.Range('E11').Value = '  +0.02%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'4.25"
$ws.Range('E12').Value = '  -0.01%  '

# Row 13
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.871.45'
$ws.Range('E13').Value = '  +0.74%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.643.16'
$ws.Range('E14').Value = '  +0.13%  '

# Row 15
$ws.Range('D15').Value = "'0.552"
$ws.Range('E15').Value = '  -1.18%  '

# Row 16
$ws.Range('D16').Value = '0.0₃0762'
$ws.Range('E16').Value = '  -0.24%  '

# Row 17
$ws.Range('D17').Value = "'63.18"
$ws.Range('E17').Value = '  +0.65%  '

# Row 18
$ws.Range('D18').Value = '26.038.97'
$ws.Range('E18').Value = '  +1.05%  '

# Row 19
$ws.Range('E19').Value = '  -0.89%  '

# Row 20
$ws.Range('D20').Value = "'4.45"
$ws.Range('E20').Value = '  +0.01%  '

# Row 21
$ws.Range('D21').Value = "'192.53"
$ws.Range('E21').Value = '  -0.87%  '

# Row 22
$ws.Range('D22').Value = "'9.99"
$ws.Range('E22').Value = '  +0.51%  '

# Row 23
$ws.Range('D23').Value = "'6.37"
$ws.Range('E23').Value = '  +1.70%  '

# Row 24
$ws.Range('D24').Value = "'0.994"
$ws.Range('E24').Value = '  -0.76%  '

# Row 25
$ws.Range('E25').Value = '  -1.72%  '

# Row 26
$ws.Range('D26').Value = "'141.69"
$ws.Range('E26').Value = '  -0.68%  '

# Row 27
$ws.Range('D27').Value = "'0.123"
$ws.Range('E27').Value = '  +0.13%  '

# Row 28
$ws.Range('D28').Value = "'6.87"
$ws.Range('E28').Value = '  -0.02%  '

# Row 29
$ws.Range('D29').Value = "'15.58"
$ws.Range('E29').Value = '  +0.18%  '

# Row 30
$ws.Range('E30').Value = '  -0.10%  '

# Row 31
$ws.Range('D31').Value = "'0.0493"
$ws.Range('E31').Value = '  -0.19%  '

# Row 32
$ws.Range('D32').Value = "'3.33"
$ws.Range('E32').Value = '  -0.02%  '

# Row 33
$ws.Range('D33').Value = "'3.24"
$ws.Range('E33').Value = '  -0.27%  '

# Row 34
$ws.Range('D34').Value = "'1.60"
$ws.Range('E34').Value = '  +1.41%  '

# Row 35
$ws.Range('E35').Value = '  -0.24%  '

# Row 36
$ws.Range('D36').Value = "'0.907"
$ws.Range('E36').Value = '  +0.62%  '

# Row 37
$ws.Range('D37').Value = '1.149.99'
$ws.Range('E37').Value = '  +1.80%  '

# Row 38
$ws.Range('D38').Value = "'0.545"
$ws.Range('E38').Value = '  -0.25%  '

# Row 39
$ws.Range('E39').Value = '  -1.72%  '

# Row 40
$ws.Range('D40').Value = "'0.0156"
$ws.Range('E40').Value = '  +0.45%  '

# Row 41
$ws.Range('E41').Value = '  -0.85%  '

# Row 42
$ws.Range('D42').Value = "'5.61"
$ws.Range('E42').Value = '  +0.67%  '

# Row 43
$ws.Range('D43').Value = "'100.28"
$ws.Range('E43').Value = '  +0.27%  '

# Row 44
$ws.Range('D44').Value = "'0.800"
$ws.Range('E44').Value = '  -0.52%  '

# Row 45
$ws.Range('D45').Value = '1.780.50'
$ws.Range('E45').Value = '  +0.73%  '

# Row 46
$ws.Range('E46').Value = '  -2.72%  '

# Row 47
$ws.Range('D47').Value = "'55.59"
$ws.Range('E47').Value = '  +0.95%  '

# Row 48
$ws.Range('D48').Value = "'0.0518"
$ws.Range('E48').Value = '  +3.07%  '

# Row 49
$ws.Range('E49').Value = '  +6.00%  '

# Row 50
$ws.Range('D50').Value = "'0.415"
$ws.Range('E50').Value = '  -0.32%  '

# Row 51
$ws.Range('D51').Value = "'7.62"
$ws.Range('E51').Value = '  +0.53%  '
